# Performance.xlsx update:
#  - "embedded array pool into transposition table"
#  - "added transposition table to perft definitely"
#
# This appends one more perft benchmark block (rows 105-107) to the
# "Initial Position Single Thread" sheet, re-using the formatting of the
# most recent existing block (rows 96-98) as a template, tags it with a
# new commit note ("added transposition table definitely"), and fixes up
# the now-stale P100 note style so it matches the rest of the column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Normalise the P100 note cell style (it drifts from the rest of the
#    "commit note" column, like P92/P96/P97 - copy their formatting).
# ---------------------------------------------------------------------
$ws1.Range("P96").Copy()
$ws1.Range("P100").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 2) Add the new perft block: rows 105-107 (row 104 stays blank, same as
#    the existing gap rows between blocks, e.g. row 99).
#    Formatting is cloned from the previous block (rows 96-98), which
#    uses identical column styling other than the G/M ratio columns.
# ---------------------------------------------------------------------
$ws1.Range("A96:N96").Copy()
$ws1.Range("A105:N105").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("C7:N7").Copy()
$ws1.Range("C106:N106").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("I7:N7").Copy()
$ws1.Range("I107:N107").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# G/M ratio columns in the template block use style 24; the new block
# uses style 35 (same as G31/M105), so touch those up explicitly.
$ws1.Range("G31").Copy()
$ws1.Range("G105").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws1.Range("G31").Copy()
$ws1.Range("M105").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 105 (depth 4 pair)
$ws1.Range("A105").Value = 46049
$ws1.Range("C105").Value = 4
$ws1.Range("D105").Value = 206603
$ws1.Range("E105").Value = 169
$ws1.Range("F105").Formula = "=D105/E105*1000"
$ws1.Range("G105").Formula = "=(E96-E105)/E96"
$ws1.Range("H105").Formula = "=(F105-80000000)/80000000"
$ws1.Range("I105").Value = 4
$ws1.Range("J105").Value = 197281
$ws1.Range("K105").Value = 7
$ws1.Range("L105").Formula = "=J105/K105*1000"
$ws1.Range("M105").Formula = "=(K96-K105)/K96"
$ws1.Range("N105").Formula = "=(L105-80000000)/80000000"
$ws1.Range("P105").Value = "added transposition table definitely"

# Row 106 (depth 5 pair)
$ws1.Range("C106").Value = 5
$ws1.Range("D106").Value = 5072212
$ws1.Range("E106").Value = 3708
$ws1.Range("F106").Formula = "=D106/E106*1000"
$ws1.Range("G106").Formula = "=(E97-E106)/E97"
$ws1.Range("H106").Formula = "=(F106-80000000)/80000000"
$ws1.Range("I106").Value = 5
$ws1.Range("J106").Value = 4880523
$ws1.Range("K106").Value = 121
$ws1.Range("L106").Formula = "=J106/K106*1000"
$ws1.Range("M106").Formula = "=(K97-K106)/K97"
$ws1.Range("N106").Formula = "=(L106-80000000)/80000000"

# Row 107 (depth 6, right-hand table only)
$ws1.Range("I107").Value = 6
$ws1.Range("J107").Value = 119060324
$ws1.Range("K107").Value = 2626
$ws1.Range("L107").Formula = "=J107/K107*1000"
$ws1.Range("M107").Formula = "=(K98-K107)/K98"
$ws1.Range("N107").Formula = "=(L107-80000000)/80000000"

# ---------------------------------------------------------------------
# 3) Move the visible selection to follow the newly added rows, as in
#    the authored workbook (selection parks on M109 after the edit).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("M109").Select()

# ---------------------------------------------------------------------
# 4) Re-touch the "vs other Enignes" sheet's note cells so their shared
#    text is rewritten after the new string is inserted (content itself
#    - stockfish livello 1 / vinto / https://lichess.org/ - is unchanged).
# ---------------------------------------------------------------------
$ws2.Range("B1").Value = "stockfish livello 1"
$ws2.Range("C1").Value = "vinto"
$ws2.Range("D1").Value = "https://lichess.org/"
